$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Episodes")

# Insert a new column before the last column (episode_tags), shifting it right
$ws.Range("AB1:AB5").Insert(-4161)

# New column header + values for "continuity_of_support"
$ws.Range("AB1").Value = "continuity_of_support"
$ws.Range("AB2").Value = 2
$ws.Range("AB3").Value = 2
$ws.Range("AB4").Value = 1
$ws.Range("AB5").Value = 9

# Make Episodes the active/selected sheet, with AC5 selected
$ws.Activate()
$ws.Range("AC5").Select()
